# Update the Efna3-Epha5 LR-pairs sheet with a refreshed TPM-based run.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the various detection-rate / expression / specificity metrics.
#
# Row 2 (ECs -> MuSCs) and row 3 (FAPs -> MuSCs, originally row 3 in the old sheet)
# get new numeric values and a new ligand/receptor column order (Efna3, Epha5 now
# sit in B/C, with the target-cluster cell holding MuSCs/Resolving-Mac in D).
# New rows 3, 5, 6 and 7 are added for the additional target cluster "Resolving-Mac"
# plus the extra "MuSCs -> MuSCs/Resolving-Mac" sending-cluster pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Efna3/Epha5 -> MuSCs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna3"
$ws.Cells.Item(2, 3).Value = "Epha5"
$ws.Cells.Item(2, 4).Value = "MuSCs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2420556666666667
$ws.Cells.Item(2, 8).Value = 0.726167
$ws.Cells.Item(2, 9).Value = 0.5314769098578004
$ws.Cells.Item(2, 10).Value = 0.5314769098578004
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.02440666666666666
$ws.Cells.Item(2, 14).Value = 0.07321999999999999
$ws.Cells.Item(2, 15).Value = 0.9884709884709885
$ws.Cells.Item(2, 16).Value = 0.9884709884709886
$ws.Cells.Item(2, 17).Value = 0.00590777197111111
$ws.Cells.Item(2, 18).Value = 0.05316994773999999
$ws.Cells.Item(2, 19).Value = 0.5253495064366464
$ws.Cells.Item(2, 20).Value = 0.5253495064366465

# Row 3: ECs -> Efna3/Epha5 -> Resolving-Mac (new row)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna3"
$ws.Cells.Item(3, 3).Value = "Epha5"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2420556666666667
$ws.Cells.Item(3, 8).Value = 0.726167
$ws.Cells.Item(3, 9).Value = 0.5314769098578004
$ws.Cells.Item(3, 10).Value = 0.5314769098578004
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.0002846666666666667
$ws.Cells.Item(3, 14).Value = 0.000854
$ws.Cells.Item(3, 15).Value = 0.01152901152901153
$ws.Cells.Item(3, 16).Value = 0.01152901152901153
$ws.Cells.Item(3, 17).Value = 0.00006890517977777778
$ws.Cells.Item(3, 18).Value = 0.000620146618
$ws.Cells.Item(3, 19).Value = 0.006127403421154003
$ws.Cells.Item(3, 20).Value = 0.006127403421154003

# Row 4: FAPs -> Efna3/Epha5 -> MuSCs (was row 3 in the old sheet)
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Efna3"
$ws.Cells.Item(4, 3).Value = "Epha5"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.2054156666666667
$ws.Cells.Item(4, 8).Value = 0.616247
$ws.Cells.Item(4, 9).Value = 0.4510271759376837
$ws.Cells.Item(4, 10).Value = 0.4510271759376837
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.02440666666666666
$ws.Cells.Item(4, 14).Value = 0.07321999999999999
$ws.Cells.Item(4, 15).Value = 0.9884709884709885
$ws.Cells.Item(4, 16).Value = 0.9884709884709886
$ws.Cells.Item(4, 17).Value = 0.005013511704444444
$ws.Cells.Item(4, 18).Value = 0.04512160534
$ws.Cells.Item(4, 19).Value = 0.4458272784264006
$ws.Cells.Item(4, 20).Value = 0.4458272784264006

# Row 5: FAPs -> Efna3/Epha5 -> Resolving-Mac (new row)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna3"
$ws.Cells.Item(5, 3).Value = "Epha5"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.2054156666666667
$ws.Cells.Item(5, 8).Value = 0.616247
$ws.Cells.Item(5, 9).Value = 0.4510271759376837
$ws.Cells.Item(5, 10).Value = 0.4510271759376837
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.0002846666666666667
$ws.Cells.Item(5, 14).Value = 0.000854
$ws.Cells.Item(5, 15).Value = 0.01152901152901153
$ws.Cells.Item(5, 16).Value = 0.01152901152901153
$ws.Cells.Item(5, 17).Value = 0.00005847499311111111
$ws.Cells.Item(5, 18).Value = 0.000526274938
$ws.Cells.Item(5, 19).Value = 0.005199897511283067
$ws.Cells.Item(5, 20).Value = 0.005199897511283067

# Row 6: MuSCs -> Efna3/Epha5 -> MuSCs (new row)
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Efna3"
$ws.Cells.Item(6, 3).Value = "Epha5"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.007968333333333332
$ws.Cells.Item(6, 8).Value = 0.023905
$ws.Cells.Item(6, 9).Value = 0.01749591420451593
$ws.Cells.Item(6, 10).Value = 0.01749591420451593
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.02440666666666666
$ws.Cells.Item(6, 14).Value = 0.07321999999999999
$ws.Cells.Item(6, 15).Value = 0.9884709884709885
$ws.Cells.Item(6, 16).Value = 0.9884709884709886
$ws.Cells.Item(6, 17).Value = 0.0001944804555555555
$ws.Cells.Item(6, 18).Value = 0.0017503241
$ws.Cells.Item(6, 19).Value = 0.01729420360794147
$ws.Cells.Item(6, 20).Value = 0.01729420360794147

# Row 7: MuSCs -> Efna3/Epha5 -> Resolving-Mac (new row)
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Efna3"
$ws.Cells.Item(7, 3).Value = "Epha5"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.007968333333333332
$ws.Cells.Item(7, 8).Value = 0.023905
$ws.Cells.Item(7, 9).Value = 0.01749591420451593
$ws.Cells.Item(7, 10).Value = 0.01749591420451593
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.0002846666666666667
$ws.Cells.Item(7, 14).Value = 0.000854
$ws.Cells.Item(7, 15).Value = 0.01152901152901153
$ws.Cells.Item(7, 16).Value = 0.01152901152901153
$ws.Cells.Item(7, 17).Value = 0.000002268318888888889
$ws.Cells.Item(7, 18).Value = 0.00002041487
$ws.Cells.Item(7, 19).Value = 0.0002017105965744607
$ws.Cells.Item(7, 20).Value = 0.0002017105965744607
